$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.172.60"
$ws.Range("E2").Value = "  +2.40%  "

$ws.Range("D3").Value = "3.642.45"
$ws.Range("E3").Value = "  +3.93%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "'606.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.98%  "

$ws.Range("D6").Value = "'202.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.22%  "

$ws.Range("E7").Value = "  +1.23%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  +10.14%  "

$ws.Range("E10").Value = "  +0.95%  "

$ws.Range("D11").Value = "'54.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.96%  "

$ws.Range("E12").Value = "  +2.75%  "

$ws.Range("D13").Value = "'9.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.08%  "

$ws.Range("D14").Value = "4.224.87"
$ws.Range("E14").Value = "  +4.12%  "

$ws.Range("D15").Value = "'680.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +13.97%  "

$ws.Range("D16").Value = "71.267.59"
$ws.Range("E16").Value = "  +2.36%  "

$ws.Range("D17").Value = "'12.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.56%  "

$ws.Range("D18").Value = "3.640.23"
$ws.Range("E18").Value = "  +3.00%  "

$ws.Range("D19").Value = "'19.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.92%  "

$ws.Range("E20").Value = "  +0.36%  "

$ws.Range("E21").Value = "  +2.25%  "

$ws.Range("D22").Value = "'18.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.06%  "

$ws.Range("D23").Value = "'5.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.03%  "

$ws.Range("D24").Value = "'105.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.21%  "

$ws.Range("D25").Value = "'4.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.74%  "

$ws.Range("D26").Value = "'3.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.35%  "

$ws.Range("D27").Value = "'10.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.71%  "

$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "'6.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.85%  "

$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'9.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.40%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'34.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.13%  "

$ws.Range("B31").Value = "dogwifhat"
$ws.Range("C31").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D31").Value = "'4.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.33%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'7.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.83%  "

$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "'12.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.116"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.84%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'63.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.66%  "

$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0878"
$ws.Range("E36").Value = "  +7.81%  "

$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "3.931.60"
$ws.Range("E37").Value = "  +4.93%  "

$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "'523.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.11%  "

$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.06%  "

$ws.Range("D40").Value = "'3.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.79%  "

$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "'3.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.24%  "

$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "'0.392"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.74%  "

$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'36.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.57%  "

$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.139"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.14%  "

$ws.Range("D45").Value = "'3.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.47%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0461"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.68%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.56%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.141"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'8.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.75%  "

$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Value = "'1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.21%  "

$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "'0.000249"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.85%  "

